$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1583.7446
$ws.Cells.Item(15, 9).Value = 1583.7446
$ws.Cells.Item(15, 11).Value = 4751.2338
$ws.Cells.Item(15, 13).Value = -4582.2338

$ws.Cells.Item(38, 8).Value = 1064.1428
$ws.Cells.Item(38, 9).Value = 389.8
$ws.Cells.Item(38, 10).Value = 2750
$ws.Cells.Item(38, 11).Value = 1169.4
$ws.Cells.Item(38, 12).Value = 8250
$ws.Cells.Item(38, 13).Value = -797.4000000000001
$ws.Cells.Item(38, 14).Value = -8994

$ws.Cells.Item(53, 8).Value = 182.68182
$ws.Cells.Item(53, 9).Value = 142.3
$ws.Cells.Item(53, 10).Value = 216.33333
$ws.Cells.Item(53, 11).Value = 142.3
$ws.Cells.Item(53, 12).Value = 216.33333
$ws.Cells.Item(53, 13).Value = 494.7
$ws.Cells.Item(53, 14).Value = -1490.33333

$ws.Cells.Item(58, 8).Value = 1792.6875
$ws.Cells.Item(58, 9).Value = 213.57143
$ws.Cells.Item(58, 10).Value = 3020.889
$ws.Cells.Item(58, 11).Value = 640.71429
$ws.Cells.Item(58, 12).Value = 9062.667000000001
$ws.Cells.Item(58, 13).Value = -490.71429
$ws.Cells.Item(58, 14).Value = -9362.667000000001

$ws.Cells.Item(111, 8).Value = 7512.375
$ws.Cells.Item(111, 9).Value = 9519.799999999999
$ws.Cells.Item(111, 10).Value = 4166.6665
$ws.Cells.Item(111, 11).Value = 28559.4
$ws.Cells.Item(111, 12).Value = 12499.9995
$ws.Cells.Item(111, 13).Value = -25492.4
$ws.Cells.Item(111, 14).Value = -18633.9995

$ws.Cells.Item(129, 8).Value = 830.56757
$ws.Cells.Item(129, 10).Value = 931.74194
$ws.Cells.Item(129, 12).Value = 2795.22582
$ws.Cells.Item(129, 14).Value = -12795.22582

$ws.Cells.Item(132, 8).Value = 2048.2122
$ws.Cells.Item(132, 9).Value = 2460.9583
$ws.Cells.Item(132, 10).Value = 947.55554
$ws.Cells.Item(132, 11).Value = 7382.874899999999
$ws.Cells.Item(132, 12).Value = 2842.66662
$ws.Cells.Item(132, 13).Value = -4852.874899999999
$ws.Cells.Item(132, 14).Value = -7902.66662

$ws.Cells.Item(137, 8).Value = 2175.8
$ws.Cells.Item(137, 9).Value = 1681.4688
$ws.Cells.Item(137, 10).Value = 3392.6155
$ws.Cells.Item(137, 11).Value = 5044.4064
$ws.Cells.Item(137, 12).Value = 10177.8465
$ws.Cells.Item(137, 13).Value = -2494.4064
$ws.Cells.Item(137, 14).Value = -15277.8465

$ws.Cells.Item(141, 8).Value = 2404.0688
$ws.Cells.Item(141, 9).Value = 1305.0454
$ws.Cells.Item(141, 10).Value = 5858.143
$ws.Cells.Item(141, 11).Value = 3915.1362
$ws.Cells.Item(141, 12).Value = 17574.429
$ws.Cells.Item(141, 13).Value = 1264.8638
$ws.Cells.Item(141, 14).Value = -27934.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(25, 8).Value = 1987.5
$ws.Cells.Item(25, 9).Value = 1650
$ws.Cells.Item(25, 10).Value = 3000
$ws.Cells.Item(25, 11).Value = 1650
$ws.Cells.Item(25, 12).Value = 3000
$ws.Cells.Item(25, 13).Value = -1248
$ws.Cells.Item(25, 14).Value = -3804

$ws.Cells.Item(110, 8).Value = 1701.5
$ws.Cells.Item(110, 9).Value = 1701.5
$ws.Cells.Item(110, 11).Value = 1701.5
$ws.Cells.Item(110, 13).Value = 343.5

$ws.Cells.Item(118, 8).Value = 31079.8
$ws.Cells.Item(118, 10).Value = 31079.8
$ws.Cells.Item(118, 12).Value = 31079.8
$ws.Cells.Item(118, 14).Value = -34393.8

$ws.Cells.Item(122, 8).Value = 3125824
$ws.Cells.Item(122, 9).Value = 655.9429
$ws.Cells.Item(122, 10).Value = 25002000
$ws.Cells.Item(122, 11).Value = 1967.8287
$ws.Cells.Item(122, 12).Value = 75006000
$ws.Cells.Item(122, 13).Value = 482.1713
$ws.Cells.Item(122, 14).Value = -75010900

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 14).ClearContents()  # N29 removed

$ws.Cells.Item(55, 8).Value = 61470
$ws.Cells.Item(55, 10).Value = 61470
$ws.Cells.Item(55, 12).Value = 61470
$ws.Cells.Item(55, 14).Value = -62016

$ws.Cells.Item(80, 8).Value = 204
$ws.Cells.Item(80, 9).Value = 181.66667
$ws.Cells.Item(80, 10).Value = 206.91304
$ws.Cells.Item(80, 11).Value = 181.66667
$ws.Cells.Item(80, 12).Value = 206.91304
$ws.Cells.Item(80, 13).Value = 816.3333299999999
$ws.Cells.Item(80, 14).Value = -2202.91304

$ws.Cells.Item(83, 8).Value = 204
$ws.Cells.Item(83, 9).Value = 181.66667
$ws.Cells.Item(83, 10).Value = 206.91304
$ws.Cells.Item(83, 11).Value = 908.3333500000001
$ws.Cells.Item(83, 12).Value = 1034.5652
$ws.Cells.Item(83, 13).Value = 4083.66665
$ws.Cells.Item(83, 14).Value = -11018.5652

$ws.Cells.Item(94, 8).Value = 1186.4445
$ws.Cells.Item(94, 9).Value = 1055.3636
$ws.Cells.Item(94, 10).Value = 1392.4286
$ws.Cells.Item(94, 11).Value = 1055.3636
$ws.Cells.Item(94, 12).Value = 1392.4286
$ws.Cells.Item(94, 13).Value = -604.3635999999999
$ws.Cells.Item(94, 14).Value = -2294.4286

$ws.Cells.Item(99, 8).Value = 1194.9032
$ws.Cells.Item(99, 9).Value = 974.1818
$ws.Cells.Item(99, 11).Value = 974.1818
$ws.Cells.Item(99, 13).Value = 523.8182

$ws.Cells.Item(134, 8).Value = 20272.145
$ws.Cells.Item(134, 9).Value = 2028.5128
$ws.Cells.Item(134, 10).Value = 64741
$ws.Cells.Item(134, 11).Value = 6085.538399999999
$ws.Cells.Item(134, 12).Value = 194223
$ws.Cells.Item(134, 13).Value = -3550.538399999999
$ws.Cells.Item(134, 14).Value = -199293

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1069.76
$ws.Cells.Item(16, 9).Value = 472
$ws.Cells.Item(16, 10).Value = 1717.3334
$ws.Cells.Item(16, 11).Value = 472
$ws.Cells.Item(16, 12).Value = 1717.3334
$ws.Cells.Item(16, 13).Value = -185
$ws.Cells.Item(16, 14).Value = -2291.3334

$ws.Cells.Item(113, 8).Value = 1069.76
$ws.Cells.Item(113, 9).Value = 472
$ws.Cells.Item(113, 10).Value = 1717.3334
$ws.Cells.Item(113, 11).Value = 472
$ws.Cells.Item(113, 12).Value = 1717.3334
$ws.Cells.Item(113, 13).Value = 1698
$ws.Cells.Item(113, 14).Value = -6057.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 696.2023
$ws.Cells.Item(113, 9).Value = 698.0328
$ws.Cells.Item(113, 10).Value = 692.2143
$ws.Cells.Item(113, 11).Value = 2094.0984
$ws.Cells.Item(113, 12).Value = 2076.6429
$ws.Cells.Item(113, 13).Value = 75.90160000000014
$ws.Cells.Item(113, 14).Value = -6416.6429

$ws.Cells.Item(120, 8).Value = 7007.8
$ws.Cells.Item(120, 9).Value = 11515
$ws.Cells.Item(120, 10).Value = 5881
$ws.Cells.Item(120, 11).Value = 34545
$ws.Cells.Item(120, 12).Value = 17643
$ws.Cells.Item(120, 13).Value = -29707
$ws.Cells.Item(120, 14).Value = -27319

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 238.07143
$ws.Cells.Item(107, 9).Value = 180.27272
$ws.Cells.Item(107, 10).Value = 450
$ws.Cells.Item(107, 11).Value = 180.27272
$ws.Cells.Item(107, 12).Value = 450
$ws.Cells.Item(107, 13).Value = 1739.72728
$ws.Cells.Item(107, 14).Value = -4290

$ws.Cells.Item(113, 8).Value = 1911.5938
$ws.Cells.Item(113, 10).Value = 2019.8235
$ws.Cells.Item(113, 12).Value = 2019.8235
$ws.Cells.Item(113, 14).Value = -6359.8235

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(21, 8).Value = 10000
$ws.Cells.Item(21, 10).Value = 10000
$ws.Cells.Item(21, 12).Value = 10000
$ws.Cells.Item(21, 14).Value = -10348

$ws.Cells.Item(61, 8).Value = 921906.2
$ws.Cells.Item(61, 9).Value = 15673.5
$ws.Cells.Item(61, 10).Value = 3338526.8
$ws.Cells.Item(61, 11).Value = 15673.5
$ws.Cells.Item(61, 12).Value = 3338526.8
$ws.Cells.Item(61, 13).Value = -15471.5
$ws.Cells.Item(61, 14).Value = -3338930.8

$ws.Cells.Item(64, 8).Value = 31570
$ws.Cells.Item(64, 10).Value = 31570
$ws.Cells.Item(64, 12).Value = 31570
$ws.Cells.Item(64, 14).Value = -32020

$ws.Cells.Item(67, 8).Value = 31570
$ws.Cells.Item(67, 10).Value = 31570
$ws.Cells.Item(67, 12).Value = 31570
$ws.Cells.Item(67, 14).Value = -33130

$ws.Cells.Item(75, 8).Value = 32000
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 32000
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = 32000
$ws.Cells.Item(75, 13).ClearContents()  # M75 removed
$ws.Cells.Item(75, 14).Value = -33872

$ws.Cells.Item(78, 8).Value = 32000
$ws.Cells.Item(78, 9).Value = 0
$ws.Cells.Item(78, 10).Value = 32000
$ws.Cells.Item(78, 11).Value = 0
$ws.Cells.Item(78, 12).Value = 96000
$ws.Cells.Item(78, 13).ClearContents()  # M78 removed
$ws.Cells.Item(78, 14).Value = -105360

$ws.Cells.Item(92, 8).Value = 18896.334
$ws.Cells.Item(92, 10).Value = 18896.334
$ws.Cells.Item(92, 12).Value = 18896.334
$ws.Cells.Item(92, 14).Value = -23888.334

$ws.Cells.Item(113, 8).Value = 921906.2
$ws.Cells.Item(113, 9).Value = 15673.5
$ws.Cells.Item(113, 10).Value = 3338526.8
$ws.Cells.Item(113, 11).Value = 15673.5
$ws.Cells.Item(113, 12).Value = 3338526.8
$ws.Cells.Item(113, 13).Value = -13503.5
$ws.Cells.Item(113, 14).Value = -3342866.8

$ws.Cells.Item(119, 8).Value = 79800
$ws.Cells.Item(119, 10).Value = 79800
$ws.Cells.Item(119, 12).Value = 79800
$ws.Cells.Item(119, 14).Value = -89476

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(33, 8).Value = 21400
$ws.Cells.Item(33, 10).Value = 21400
$ws.Cells.Item(33, 12).Value = 21400
$ws.Cells.Item(33, 14).Value = -21900

$ws.Cells.Item(36, 8).Value = 21400
$ws.Cells.Item(36, 10).Value = 21400
$ws.Cells.Item(36, 12).Value = 21400
$ws.Cells.Item(36, 14).Value = -21900

$ws.Cells.Item(70, 8).Value = 37303.89
$ws.Cells.Item(70, 10).Value = 37303.89
$ws.Cells.Item(70, 12).Value = 37303.89
$ws.Cells.Item(70, 14).Value = -37933.89

$ws.Cells.Item(73, 8).Value = 37303.89
$ws.Cells.Item(73, 10).Value = 37303.89
$ws.Cells.Item(73, 12).Value = 37303.89
$ws.Cells.Item(73, 14).Value = -39487.89

$ws.Cells.Item(126, 8).Value = 1041.7667
$ws.Cells.Item(126, 9).Value = 981.375
$ws.Cells.Item(126, 10).Value = 1283.3334
$ws.Cells.Item(126, 11).Value = 2944.125
$ws.Cells.Item(126, 12).Value = 3850.0002
$ws.Cells.Item(126, 13).Value = -474.125
$ws.Cells.Item(126, 14).Value = -8790.0002

$ws.Cells.Item(132, 8).Value = 1425.3043
$ws.Cells.Item(132, 9).Value = 741.0968
$ws.Cells.Item(132, 11).Value = 2223.2904
$ws.Cells.Item(132, 13).Value = 306.7096000000001
